$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a number to Excel's auto-detection;
# force them to Text format first so they are stored as literal strings,
# matching the inlineStr cells in the original workbook, then restore the
# default "Normal" style so the saved style index is unchanged.
$numericLooking = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '20.489.76'
$ws.Range('E2').Value = '  +2.60%  '
$ws.Range('D3').Value = '1.472.16'
$ws.Range('E3').Value = '  +3.81%  '
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').Value = '0.9738'
$ws.Range('E5').Value = '  -2.67%  '
$ws.Range('D6').Value = '275.10'
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('D7').Value = '0.3654'
$ws.Range('E7').Value = '  -0.94%  '
$ws.Range('D8').Value = '0.3069'
$ws.Range('E8').Value = '  -1.15%  '
$ws.Range('D9').Value = '39.87'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').Value = '1.050'
$ws.Range('E10').Value = '  +0.77%  '
$ws.Range('D11').Value = '0.06612'
$ws.Range('E11').Value = '  +1.06%  '
$ws.Range('D12').Value = '0.9973'
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('D13').Value = '18.13'
$ws.Range('E13').Value = '  +2.45%  '
$ws.Range('D14').Value = '5.458'
$ws.Range('E14').Value = '  -0.88%  '
$ws.Range('D15').Value = '6.156'
$ws.Range('E15').Value = '  -1.02%  '
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').Value = '1.473.27'
$ws.Range('E17').Value = '  +3.62%  '
$ws.Range('D18').Value = '0.9836'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').Value = '0.05900'
$ws.Range('E19').Value = '  +3.65%  '
$ws.Range('D20').Value = '69.67'
$ws.Range('E20').Value = '  -2.32%  '
$ws.Range('D21').Value = '5.459'
$ws.Range('E21').Value = '  -2.89%  '
$ws.Range('D22').Value = '14.41'
$ws.Range('E22').Value = '  -2.55%  '
$ws.Range('D23').Value = '10.93'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('D24').Value = '2.249'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').Value = '20.505.13'
$ws.Range('E25').Value = '  +2.46%  '
$ws.Range('D26').Value = '142.43'
$ws.Range('E26').Value = '  +6.94%  '
$ws.Range('D27').Value = '2.134'
$ws.Range('E27').Value = '  -6.66%  '
$ws.Range('D28').Value = '17.29'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').Value = '1.626.09'
$ws.Range('E29').Value = '  +2.82%  '
$ws.Range('D30').Value = '114.07'
$ws.Range('E30').Value = '  +3.52%  '
$ws.Range('D31').Value = '3.853'
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('D32').Value = '4.976'
$ws.Range('E32').Value = '  -5.09%  '
$ws.Range('D33').Value = '0.8018'
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').Value = '0.07858'
$ws.Range('E34').Value = '  +1.19%  '
$ws.Range('D35').Value = '1.536'
$ws.Range('E35').Value = '  +3.47%  '
$ws.Range('D36').Value = '0.05770'
$ws.Range('E36').Value = '  -1.70%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').Value = '1.154'
$ws.Range('E37').Value = '  +4.74%  '
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '4.728'
$ws.Range('E38').Value = '  -3.59%  '
$ws.Range('D39').Value = '7.770'
$ws.Range('E39').Value = '  -5.23%  '
$ws.Range('D40').Value = '0.9752'
$ws.Range('E40').Value = '  -2.46%  '
$ws.Range('D41').Value = '0.02041'
$ws.Range('E41').Value = '  -0.79%  '
$ws.Range('D42').Value = '10.42'
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('D43').Value = '0.1871'
$ws.Range('E43').Value = '  -0.66%  '
$ws.Range('D44').Value = '0.5287'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').Value = '3.489'
$ws.Range('E45').Value = '  -1.40%  '
$ws.Range('D46').Value = '12.03'
$ws.Range('E46').Value = '  -3.04%  '
$ws.Range('D47').Value = '117.75'
$ws.Range('E47').Value = '  +0.68%  '
$ws.Range('D48').Value = '0.5191'
$ws.Range('E48').Value = '  -0.37%  '
$ws.Range('D49').Value = '1.772'
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').Value = '0.06453'
$ws.Range('E50').Value = '  +4.15%  '
$ws.Range('D51').Value = '0.9876'
$ws.Range('E51').Value = '  -1.31%  '

foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}
